# Weekly update: a new record (week of 2022-12-12) is inserted at the top
# of the data block (row 50), pushing every existing record at row 50
# and below down by one row. The previously-last record (old row 92)
# becomes the new last record (row 93).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50; Excel shifts rows 50:92 down to 51:93
# and extends the used range / dimension automatically.
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new weekly record.
$ws.Range("A50").Value() = 7
$ws.Range("B50").Value() = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C50").Value() = "Ñuble"
$ws.Range("D50").Value() = 44907
$ws.Range("E50").Value() = 16
$ws.Range("F50").Value() = 100112022
$ws.Range("G50").Value() = "Arveja Verde"
$ws.Range("H50").Value() = "Sin especificar"
$ws.Range("I50").Value() = "Primera"
$ws.Range("J50").Value() = 60
$ws.Range("K50").Value() = 22000
$ws.Range("L50").Value() = 23000
$ws.Range("M50").Value() = 22500
$ws.Range("N50").Value() = "$/saco 25 kilos"
$ws.Range("O50").Value() = "Región del Maule"
$ws.Range("P50").Value() = 900
$ws.Range("Q50").Value() = 25
$ws.Range("R50").Value() = "Hortaliza"
